# CRM-1761 Add customer mobile no in pending booking and pending spare on
# partner panel.
#
# The "Spare Requested Parts" export template gets a new "Customer Phone
# Number" column inserted right after "Customer Name" (i.e. at column B),
# pushing every following column one slot to the right. Row 1 holds the
# human-readable header, row 2 holds the merge-field placeholder that the
# backend substitutes when generating the real spreadsheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column B ("Invoice Id"). The new
# column inherits column A's look (bold, centered, bordered header / plain
# data row), which is exactly the visual style already used by the sheet's
# other "phone number" column.
$ws.Columns("B").Insert()

# Header row.
$ws.Range("B1").Value = "Customer Phone Number"

# Placeholder row consumed by the report generator.
$ws.Range("B2").Value = "{spare:customer_mobile}"

# Match the new column's width to column A so the header isn't clipped any
# worse than the template's other columns.
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth
